$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Status" column header to "Resolution"
$ws.Range("E1").Value = "Resolution"

# Replace the status values: Done -> Fixed, Todo -> Unresolved
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Value = "Fixed"
}
for ($r = 10; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = "Unresolved"
}

# Update column E width to fit the new content (closest achievable value
# to the authored 10.5703125 given this engine's width quantization)
$ws.Columns.Item(5).ColumnWidth = 9.6

# Update the active selection to match the authored change
$ws.Range("E18").Select() | Out-Null
